$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Strip stray leading spaces from a handful of Kannada / Hindi / Tamil
# field-name translations (lang_code rows for kan / hin / tam) that were
# accidentally prefixed with a space character. Several of the corrected
# strings collapse onto already-existing shared strings elsewhere in the
# sheet, which Excel automatically de-duplicates when the workbook is saved.
# ---------------------------------------------------------------------------

# kan (Kannada) section
$ws.Range("E123").Value = 'ಪ್ರದೇಶ'
$ws.Range("E124").Value = 'ಪ್ರಾಂತ್ಯ'
$ws.Range("E125").Value = 'ನಗರ'
$ws.Range("E142").Value = 'ಪ್ರಾಂತ್ಯ'
$ws.Range("E143").Value = 'ನಗರ'
$ws.Range("E157").Value = 'ಕೋಡ್'
$ws.Range("E158").Value = 'ಕೋಡ್'
$ws.Range("E159").Value = 'ಕೋಡ್'
$ws.Range("E160").Value = 'ಕೋಡ್'
$ws.Range("E161").Value = 'ಕೋಡ್'

# hin (Hindi) section
$ws.Range("E162").Value = 'देश'
$ws.Range("E163").Value = 'क्षेत्र'
$ws.Range("E165").Value = 'शहर'
$ws.Range("E166").Value = 'क्षेत्र'
$ws.Range("E197").Value = 'डाक कोड'
$ws.Range("E198").Value = 'डाक कोड'
$ws.Range("E199").Value = 'डाक कोड'
$ws.Range("E200").Value = 'डाक कोड'
$ws.Range("E201").Value = 'डाक कोड'

# tam (Tamil) section
$ws.Range("E202").Value = 'நாடு'

# ---------------------------------------------------------------------------
# Restore the default view flags (these round-trip to the same values the
# workbook already had) and move the active selection / scroll position to
# reflect where the editor was last working in the sheet.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.DisplayFormulas = $false
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true
$excel.ActiveWindow.DisplayZeros = $true
$excel.ActiveWindow.DisplayOutline = $true

$ws.Range("E126").Select()
$excel.ActiveWindow.ScrollRow = 230
$excel.ActiveWindow.ScrollColumn = 1
